$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to update (Price column D and Volume(1h) column E), taken from the
# refreshed cryptocurrency data feed.
$updates = @{
    "D2" = "63.867.79"
    "E2" = "  -0.92%  "
    "D3" = "3.054.02"
    "E3" = "  -1.45%  "
    "E4" = "  -0.06%  "
    "D5" = "559.89"
    "E5" = "  +0.25%  "
    "D6" = "142.72"
    "E6" = "  -0.81%  "
    "E7" = "  +0.06%  "
    "D8" = "3.051.58"
    "E8" = "  -1.42%  "
    "D10" = "0.153"
    "E10" = "  +1.12%  "
    "D11" = "6.28"
    "E11" = "  -10.54%  "
    "D12" = "0.490"
    "E12" = "  +6.49%  "
    "D13" = "0.0000230"
    "E13" = "  +1.62%  "
    "D14" = "35.75"
    "E14" = "  +1.90%  "
    "D15" = "3.553.54"
    "E15" = "  -1.23%  "
    "D16" = "63.937.30"
    "E16" = "  -0.94%  "
    "D17" = "3.055.43"
    "E17" = "  -1.31%  "
    "D18" = "0.110"
    "E18" = "  +0.75%  "
    "E19" = "  +1.12%  "
    "D20" = "476.50"
    "E20" = "  -1.64%  "
    "D21" = "14.07"
    "E21" = "  +2.26%  "
    "D22" = "14.77"
    "E22" = "  +11.78%  "
    "E23" = "  +1.91%  "
    "E24" = "  -0.30%  "
    "D25" = "82.72"
    "E25" = "  +2.26%  "
    "E26" = "  +0.01%  "
    "E27" = "  -0.55%  "
    "D28" = "8.15"
    "E28" = "  +1.98%  "
    "E29" = "  -1.23%  "
    "D30" = "0.999"
    "E30" = "  -0.17%  "
    "D31" = "26.32"
    "E31" = "  +0.94%  "
    "E32" = "  -0.57%  "
    "E33" = "  +0.61%  "
    "E34" = "  +0.14%  "
    "E35" = "  +1.81%  "
    "E36" = "  -1.18%  "
    "D37" = "0.0411"
    "E37" = "  +0.89%  "
    "D38" = "447.92"
    "E38" = "  -3.32%  "
    "E39" = "  -1.23%  "
    "D40" = "2.81"
    "E40" = "  +4.35%  "
    "D41" = "3.022.77"
    "E41" = "  -0.11%  "
    "E42" = "  +0.76%  "
    "D43" = "8.28"
    "E43" = "  -0.18%  "
    "E44" = "  +4.24%  "
    "D45" = "28.36"
    "E45" = "  +0.93%  "
    "D46" = "2.27"
    "E46" = "  +8.61%  "
    "E48" = "  +0.87%  "
    "D49" = "117.92"
    "E49" = "  -0.88%  "
    "E50" = "  -0.48%  "
    "E51" = "  +1.29%  "
}

# These Price cells contain values that look like plain decimal numbers
# (single dot, e.g. "559.89"), so Excel would otherwise auto-convert them to
# numeric cells. Force them to remain text, matching the rest of the column
# (other Price cells use "." as a thousands separator, e.g. "63.867.79",
# which Excel cannot interpret as a number anyway and keeps as text).
$textCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D18", "D20", "D21", "D22", "D25", "D28", "D30", "D31", "D37", "D38", "D40", "D43", "D45", "D46", "D49")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
